{"js": "const replacements = [\n  [\"2025-03-17 Monday\", \"2025-03-18 Tuesday\"],\n  [\"881\u00d78=\", \"147\u00d76=\"],\n  [\"989\u00d75=\", \"213\u00d75=\"],\n  [\"273\u00d72=\", \"130\u00d78=\"],\n  [\"894\u00d73=\", \"209\u00d77=\"],\n  [\"721\u00d75=\", \"250\u00d77=\"],\n  [\"754\u00d72=\", \"648\u00d78=\"],\n  [\"285\u00d79=\", \"221\u00d74=\"],\n  [\"720\u00d76=\", \"910\u00d75=\"],\n  [\"359\u00d77=\", \"761\u00d74=\"],\n  [\"790\u00d75=\", \"743\u00d76=\"],\n  [\"988\u00d72=\", \"887\u00d75=\"],\n  [\"131\u00d74=\", \"983\u00d72=\"],\n  [\"525\u00d74=\", \"341\u00d77=\"],\n  [\"897\u00d76=\", \"111\u00d79=\"],\n  [\"219\u00d76=\", \"922\u00d77=\"],\n  [\"532\u00d73=\", \"911\u00d72=\"],\n  [\"194\u00d74=\", \"908\u00d72=\"],\n  [\"216\u00d73=\", \"169\u00d75=\"],\n  [\"757\u00d73=\", \"919\u00d77=\"],\n  [\"925\u00d77=\", \"311\u00d72=\"],\n  [\"869\u00d78=\", \"249\u00d77=\"],\n  [\"930\u00d76=\", \"864\u00d75=\"],\n  [\"388\u00d77=\", \"751\u00d76=\"],\n  [\"274\u00d75=\", \"684\u00d78=\"],\n  [\"220\u00d74=\", \"114\u00d72=\"],\n];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n    totalReplaced++;\n  }\n  await context.sync();\n}\nreturn `replaced ${totalReplaced} of ${replacements.length}`;", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"2025-03-17 Monday\"; New=\"2025-03-18 Tuesday\"},\n    @{Old=\"881\u00d78=\"; New=\"147\u00d76=\"},\n    @{Old=\"989\u00d75=\"; New=\"213\u00d75=\"},\n    @{Old=\"273\u00d72=\"; New=\"130\u00d78=\"},\n    @{Old=\"894\u00d73=\"; New=\"209\u00d77=\"},\n    @{Old=\"721\u00d75=\"; New=\"250\u00d77=\"},\n    @{Old=\"754\u00d72=\"; New=\"648\u00d78=\"},\n    @{Old=\"285\u00d79=\"; New=\"221\u00d74=\"},\n    @{Old=\"720\u00d76=\"; New=\"910\u00d75=\"},\n    @{Old=\"359\u00d77=\"; New=\"761\u00d74=\"},\n    @{Old=\"790\u00d75=\"; New=\"743\u00d76=\"},\n    @{Old=\"988\u00d72=\"; New=\"887\u00d75=\"},\n    @{Old=\"131\u00d74=\"; New=\"983\u00d72=\"},\n    @{Old=\"525\u00d74=\"; New=\"341\u00d77=\"},\n    @{Old=\"897\u00d76=\"; New=\"111\u00d79=\"},\n    @{Old=\"219\u00d76=\"; New=\"922\u00d77=\"},\n    @{Old=\"532\u00d73=\"; New=\"911\u00d72=\"},\n    @{Old=\"194\u00d74=\"; New=\"908\u00d72=\"},\n    @{Old=\"216\u00d73=\"; New=\"169\u00d75=\"},\n    @{Old=\"757\u00d73=\"; New=\"919\u00d77=\"},\n    @{Old=\"925\u00d77=\"; New=\"311\u00d72=\"},\n    @{Old=\"869\u00d78=\"; New=\"249\u00d77=\"},\n    @{Old=\"930\u00d76=\"; New=\"864\u00d75=\"},\n    @{Old=\"388\u00d77=\"; New=\"751\u00d76=\"},\n    @{Old=\"274\u00d75=\"; New=\"684\u00d78=\"},\n    @{Old=\"220\u00d74=\"; New=\"114\u00d72=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
